$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.754495054335223
$ws.Range("C2").Value = 1.956621400508425
$ws.Range("D2").Value = 0.6868305614177075
$ws.Range("E2").Value = 0.2794402107012743
$ws.Range("G2").Value = 0.002627713015233121
$ws.Range("I2").Value = 2.448129913967307
$ws.Range("J2").Value = 0.1448514546007971
$ws.Range("N2").Value = 2.455164247716212
$ws.Range("B3").Value = 3.586508059874916
$ws.Range("C3").Value = 1.872547643677422
$ws.Range("D3").Value = 0.6788607085536285
$ws.Range("E3").Value = 0.27536416209756
$ws.Range("G3").Value = 0.002636509792782885
$ws.Range("I3").Value = 2.412306665499898
$ws.Range("J3").Value = 0.1420321755307015
$ws.Range("N3").Value = 2.455157177253682
$ws.Range("B4").Value = 3.486160846921166
$ws.Range("C4").Value = 1.822378721579923
$ws.Range("D4").Value = 0.6743822999516738
$ws.Range("E4").Value = 0.2730385314095685
$ws.Range("G4").Value = 0.002642181172548103
$ws.Range("I4").Value = 2.391667677052283
$ws.Range("J4").Value = 0.1403988741415958
$ws.Range("N4").Value = 2.455851492772439
$ws.Range("B5").Value = 3.445962712284199
$ws.Range("C5").Value = 1.802294186839561
$ws.Range("D5").Value = 0.672660964686429
$ws.Range("E5").Value = 0.2721349828750732
$ws.Range("G5").Value = 0.002644560533786353
$ws.Range("I5").Value = 2.383594615067892
$ws.Range("J5").Value = 0.1397576256559176
$ws.Range("N5").Value = 2.456308383506027
$ws.Range("B6").Value = 3.439329501659699
$ws.Range("C6").Value = 1.798980729611799
$ws.Range("D6").Value = 0.67238137764204
$ws.Range("E6").Value = 0.2719876063315141
$ws.Range("G6").Value = 0.002644959754455045
$ws.Range("I6").Value = 2.382274367892265
$ws.Range("J6").Value = 0.1396526095797057
$ws.Range("N6").Value = 2.456394705957109
$ws.Range("B7").Value = 3.485615922250133
$ws.Range("C7").Value = 1.82210640536664
$ws.Range("D7").Value = 0.6743586666858334
$ws.Range("E7").Value = 0.2730261674887018
$ws.Range("G7").Value = 0.002642212984756102
$ws.Range("I7").Value = 2.391557439172573
$ws.Range("J7").Value = 0.1403901278265494
$ws.Range("N7").Value = 2.455856952502756
$ws.Range("B8").Value = 3.69598681401277
$ws.Range("C8").Value = 1.927327834484402
$ws.Range("D8").Value = 0.6839959496269046
$ws.Range("E8").Value = 0.2779978130650633
$ws.Range("G8").Value = 0.002630690263689053
$ws.Range("I8").Value = 2.435494261042834
$ws.Range("J8").Value = 0.1438589298409738
$ws.Range("N8").Value = 2.455015751709936
$ws.Range("B9").Value = 4.131172722617976
$ws.Range("C9").Value = 2.145462353742062
$ws.Range("D9").Value = 0.7062239693675565
$ws.Range("E9").Value = 0.2891698423637195
$ws.Range("G9").Value = 0.002610223628610161
$ws.Range("I9").Value = 2.532592183312659
$ws.Range("J9").Value = 0.1514485902281919
$ws.Range("N9").Value = 2.458989834403297
$ws.Range("B10").Value = 4.465387711142057
$ws.Range("C10").Value = 2.31331444436853
$ws.Range("D10").Value = 0.7246377207229102
$ws.Range("E10").Value = 0.2982713806149704
$ws.Range("G10").Value = 0.002596465331980066
$ws.Range("I10").Value = 2.610852465893856
$ws.Range("J10").Value = 0.1575221081329659
$ws.Range("N10").Value = 2.465451996533147
$ws.Range("B11").Value = 4.620729134577687
$ws.Range("C11").Value = 2.391412293040162
$ws.Range("D11").Value = 0.7334789910518964
$ws.Range("E11").Value = 0.3026118524972787
$ws.Range("G11").Value = 0.002590479715747151
$ws.Range("I11").Value = 2.648015435700287
$ws.Range("J11").Value = 0.1603970535312698
$ws.Range("N11").Value = 2.469186361983702
$ws.Range("B12").Value = 4.680040374257374
$ws.Range("C12").Value = 2.421243477682424
$ws.Range("D12").Value = 0.7368947586526247
$ws.Range("E12").Value = 0.3042847349388538
$ws.Range("G12").Value = 0.002588252059082395
$ws.Range("I12").Value = 2.662317317986691
$ws.Range("J12").Value = 0.1615021543804005
$ws.Range("N12").Value = 2.470716848197839
$ws.Range("B13").Value = 4.667244799427067
$ws.Range("C13").Value = 2.414807246945884
$ws.Range("D13").Value = 0.7361560826138884
$ws.Range("E13").Value = 0.3039231418326551
$ws.Range("G13").Value = 0.002588730096521905
$ws.Range("I13").Value = 2.65922689136454
$ws.Range("J13").Value = 0.161263415829012
$ws.Range("N13").Value = 2.470382021582992
$ws.Range("B14").Value = 4.625598882335566
$ws.Range("C14").Value = 2.39386132741754
$ws.Range("D14").Value = 0.7337586446035402
$ws.Range("E14").Value = 0.3027488928401993
$ws.Range("G14").Value = 0.002590295665819919
$ws.Range("I14").Value = 2.649187443244358
$ws.Range("J14").Value = 0.1604876399496504
$ws.Range("N14").Value = 2.469309931806265
$ws.Range("B15").Value = 4.600153300618103
$ws.Range("C15").Value = 2.381065051960604
$ws.Range("D15").Value = 0.7322989986766117
$ws.Range("E15").Value = 0.3020334532479225
$ws.Range("G15").Value = 0.002591259687698328
$ws.Range("I15").Value = 2.643067957764345
$ws.Range("J15").Value = 0.1600146026070632
$ws.Range("N15").Value = 2.468668462467406
$ws.Range("B16").Value = 4.455303271850539
$ws.Range("C16").Value = 2.308246183293249
$ws.Range("D16").Value = 0.7240693546217756
$ws.Range("E16").Value = 0.297991786869261
$ws.Range("G16").Value = 0.002596861975536389
$ws.Range("I16").Value = 2.608455588114055
$ws.Range("J16").Value = 0.1573365047786268
$ws.Range("N16").Value = 2.465224115619947
$ws.Range("B17").Value = 4.367297130826387
$ws.Range("C17").Value = 2.264024919635176
$ws.Range("D17").Value = 0.7191404292593688
$ws.Range("E17").Value = 0.2955639383747481
$ws.Range("G17").Value = 0.002600368528823842
$ws.Range("I17").Value = 2.587625217758912
$ws.Range("J17").Value = 0.1557224973483358
$ws.Range("N17").Value = 2.463316131425785
$ws.Range("B18").Value = 4.316988732363143
$ws.Range("C18").Value = 2.238753476360785
$ws.Range("D18").Value = 0.7163491294579671
$ws.Range("E18").Value = 0.2941863135177201
$ws.Range("G18").Value = 0.002602411129532356
$ws.Range("I18").Value = 2.575790886319993
$ws.Range("J18").Value = 0.1548046928472218
$ws.Range("N18").Value = 2.462293382403885
$ws.Range("B19").Value = 4.30000814048708
$ws.Range("C19").Value = 2.230224868164328
$ws.Range("D19").Value = 0.7154115186258707
$ws.Range("E19").Value = 0.2937230901443044
$ws.Range("G19").Value = 0.002603107146256661
$ws.Range("I19").Value = 2.571809054807957
$ws.Range("J19").Value = 0.1544957393681301
$ws.Range("N19").Value = 2.461959859433406
$ws.Range("B20").Value = 4.376633326648061
$ws.Range("C20").Value = 2.268715382304777
$ws.Range("D20").Value = 0.7196605935695857
$ws.Range("E20").Value = 0.2958204370019004
$ws.Range("G20").Value = 0.002599992590017059
$ws.Range("I20").Value = 2.589827426577699
$ws.Range("J20").Value = 0.1558932192065186
$ws.Range("N20").Value = 2.463511496259571
$ws.Range("B21").Value = 4.637817992266719
$ws.Range("C21").Value = 2.400006615112375
$ws.Range("D21").Value = 0.7344609833643574
$ws.Range("E21").Value = 0.3030930009278663
$ws.Range("G21").Value = 0.0025898347652885
$ws.Range("I21").Value = 2.652130019414898
$ws.Range("J21").Value = 0.1607150559765529
$ws.Range("N21").Value = 2.469621655637241
$ws.Range("B22").Value = 4.811360383845681
$ws.Range("C22").Value = 2.487315648243907
$ws.Range("D22").Value = 0.7445293690312837
$ws.Range("E22").Value = 0.3080166988211701
$ws.Range("G22").Value = 0.002583423022048325
$ws.Range("I22").Value = 2.694185728590782
$ws.Range("J22").Value = 0.1639622701454329
$ws.Range("N22").Value = 2.474294252347391
$ws.Range("B23").Value = 4.718473505162592
$ws.Range("C23").Value = 2.440577367238689
$ws.Range("D23").Value = 0.7391191780637598
$ws.Range("E23").Value = 0.3053730574930853
$ws.Range("G23").Value = 0.002586824422464949
$ws.Range("I23").Value = 2.671615931164638
$ws.Range("J23").Value = 0.1622202957651382
$ws.Range("N23").Value = 2.47173754618268
$ws.Range("B24").Value = 4.37241153787113
$ws.Range("C24").Value = 2.26659435135673
$ws.Range("D24").Value = 0.7194252952732256
$ws.Range("E24").Value = 0.2957044173881869
$ws.Range("G24").Value = 0.002600162468878693
$ws.Range("I24").Value = 2.588831368108089
$ws.Range("J24").Value = 0.1558160043872903
$ws.Range("N24").Value = 2.463422940932389
$ws.Range("B25").Value = 4.010949843362482
$ws.Range("C25").Value = 2.085147890639576
$ws.Range("D25").Value = 0.6998488746243652
$ws.Range("E25").Value = 0.2859925499084994
$ws.Range("G25").Value = 0.002615534451425143
$ws.Range("I25").Value = 2.505127718947463
$ws.Range("J25").Value = 0.149309324880079
$ws.Range("N25").Value = 2.457301040665442
